# Insert a new data row into the "Arveja Verde" price list.
# The new record is inserted as row 63 (pushing all following rows down
# by one), which matches the reference diff: dimension grows from
# A1:R135 to A1:R136 and every row from the former 63 onward is shifted
# down by one, with a brand-new record occupying the vacated row 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 63..135 down to 64..136.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new observation.
$ws.Cells.Item(63, 1).Value  = 4
$ws.Cells.Item(63, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(63, 3).Value  = "Los Lagos"
$ws.Cells.Item(63, 4).Value  = 44810
$ws.Cells.Item(63, 5).Value  = 10
$ws.Cells.Item(63, 6).Value  = 100112022
$ws.Cells.Item(63, 7).Value  = "Arveja Verde"
$ws.Cells.Item(63, 8).Value  = "Perfection"
$ws.Cells.Item(63, 9).Value  = "Primera"
$ws.Cells.Item(63, 10).Value = 70
$ws.Cells.Item(63, 11).Value = 42000
$ws.Cells.Item(63, 12).Value = 42000
$ws.Cells.Item(63, 13).Value = 42000
$ws.Cells.Item(63, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(63, 16).Value = 1680
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
